# Delivery Roadmap Template v2 - "Release Summary" slide
# Commit: Taller rows for full summaries, added Jira link
#
# - Each "MGS-XXXXX: description..." card is split into two paragraphs:
#   a bold ticket-id line, and a smaller/lighter description line.
# - Card rows grow from 347472 -> 438912 EMU to fit the full text, and
#   every row below the first in a column shifts down to match.
# - The four (five) background "Rounded Rectangle" columns grow taller
#   to keep containing the taller stack of cards.
# - The last 3 cards in the "MAR 2027" column are removed, and the
#   "+N more" caption is reworded and moved down to its new slot.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------
# 1) Card text boxes: taller rows, bold id / lighter description
# ---------------------------------------------------------------
$cards = @(
    # Column: PHASE 2a (2)
    @{Name="TextBox 6";  Y=$null;    Id="MGS-12826"; Desc="Enhance OLBM document upload messaging for case types"}
    @{Name="TextBox 7";  Y=2267712;  Id="MGS-12838"; Desc="Phase 2a - Product fee and LTV Enhancements"}

    # Column: JUNE (5)
    @{Name="TextBox 11"; Y=$null;    Id="MGS-12833"; Desc="Enhance document upload capabilities"}
    @{Name="TextBox 12"; Y=2267712;  Id="MGS-12830"; Desc="DD Bank validation ESB Service"}
    @{Name="TextBox 13"; Y=2706624;  Id="MGS-12842"; Desc="CRA enhancements - June target"}
    @{Name="TextBox 14"; Y=3145536;  Id="MGS-13102"; Desc="Prevent Special Characters entered in OLBM causing robots to fail"}
    @{Name="TextBox 15"; Y=3584448;  Id="MGS-13079"; Desc="Remove restriction on CC seeing property valuation details"}

    # Column: SEPT (5)
    @{Name="TextBox 19"; Y=$null;    Id="MGS-13098"; Desc="Increase levels of robot validation"}
    @{Name="TextBox 20"; Y=2267712;  Id="MGS-13105"; Desc="Robot to park data entry exceptions until all acceptable data entered"}
    @{Name="TextBox 21"; Y=2706624;  Id="MGS-13077"; Desc="Remove restriction on BDM seeing property valuation details"}
    @{Name="TextBox 22"; Y=3145536;  Id="MGS-13070"; Desc="Make case status visible in Case Tracker"}
    @{Name="TextBox 23"; Y=3584448;  Id="MGS-12827"; Desc="Remove Valuation Visibility Restrictions"}

    # Column: NOV (2)
    @{Name="TextBox 27"; Y=$null;    Id="MGS-12832"; Desc="Existing Customer Online Application"}
    @{Name="TextBox 28"; Y=2267712;  Id="MGS-12831"; Desc="Worldpay to Global Pay"}

    # Column: MAR 2027 (19) - first 9 cards keep, last 3 are deleted below
    @{Name="TextBox 32"; Y=$null;    Id="MGS-13067"; Desc="Create new fields to identify Second/Holiday Home's and pass data to PCO"}
    @{Name="TextBox 33"; Y=2267712;  Id="MGS-13066"; Desc="Create Variable Income fields and pass data to PCO"}
    @{Name="TextBox 34"; Y=2706624;  Id="MGS-13062"; Desc="Create new fields for Joint Borrower Sole Owner cases and pass to PCO"}
    @{Name="TextBox 35"; Y=3145536;  Id="MGS-13063"; Desc="Create additional fields for Self Employed (Main Income)"}
    @{Name="TextBox 36"; Y=3584448;  Id="MGS-12823"; Desc="Create Secondary Income fields and pass to PCO"}
    @{Name="TextBox 37"; Y=4023360;  Id="MGS-13096"; Desc="Improve clarity of messages in post-document review emails"}
    @{Name="TextBox 38"; Y=4462272;  Id="MGS-12841"; Desc="Replacement Mortgage Illustrator tool"}
    @{Name="TextBox 39"; Y=4901184;  Id="MGS-13092"; Desc="Add date of valuation to MAP chaser letter"}
    @{Name="TextBox 40"; Y=5340096;  Id="MGS-13086"; Desc="Improve MAP email templates"}
)

foreach ($card in $cards) {
    $sh = $s.Shapes.Item($card.Name)

    if ($card.Y -ne $null) {
        $sh.Top = $card.Y / $EMU_PER_PT
    }

    $sh.TextFrame.TextRange.Text = $card.Id + "`r" + $card.Desc

    $idRange = $sh.TextFrame.TextRange.Paragraphs(1, 1)
    $idRange.Font.Bold = 1

    $descRange = $sh.TextFrame.TextRange.Paragraphs(2, 1)
    $descRange.Font.Size = 7
    $descRange.Font.Color.RGB = 8019036   # 5C5C7A

    $sh.Height = 438912 / $EMU_PER_PT
}

# ---------------------------------------------------------------
# 2) Background column rectangles grow taller
# ---------------------------------------------------------------
$containers = @("Rounded Rectangle 3", "Rounded Rectangle 8", "Rounded Rectangle 16", "Rounded Rectangle 24", "Rounded Rectangle 29")
foreach ($name in $containers) {
    $sh = $s.Shapes.Item($name)
    $sh.Height = 5303520 / $EMU_PER_PT
}

# ---------------------------------------------------------------
# 3) Drop the last three MAR 2027 cards (now summarised by "+10 more")
# ---------------------------------------------------------------
$s.Shapes.Item("TextBox 41").Delete()
$s.Shapes.Item("TextBox 42").Delete()
$s.Shapes.Item("TextBox 43").Delete()

# ---------------------------------------------------------------
# 4) "+N more" caption: reworded, resized down to 7pt, moved down
# ---------------------------------------------------------------
$more = $s.Shapes.Item("TextBox 44")
$more.Top = 6217920 / $EMU_PER_PT
$more.TextFrame.TextRange.Text = "+10 more - see Jira"
$more.TextFrame.TextRange.Font.Size = 7
